$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column I ("roboticS1Prep") currently holds the text "No" for every data row (2-41).
# Convert it into a real boolean column, formatted to display TRUE/FALSE.
$srcFont = $ws.Range("I2").Font
$dstFont = $ws.Range("I28:I41").Font
$dstFont.Name = $srcFont.Name
$dstFont.Size = $srcFont.Size
$dstFont.Bold = $srcFont.Bold
$dstFont.Italic = $srcFont.Italic
$dstFont.Color = $srcFont.Color

$rng = $ws.Range("I2:I41")
$rng.Value = $false
$rng.NumberFormat = '"TRUE";"TRUE";"FALSE"'

# Move the active selection / view to the column that was just edited.
$ws.Range("I2:I41").Select()
$excel.ActiveWindow.ScrollRow = 10
